# Generate Report for Handback
# Updates the "b3b653f4-b0ac-4d24-bf69-21d55057a0d5" row (row 6) on the
# zh-cn and de-de sheets to reflect a handback attempt that failed because
# the handback file version was stale, and records the new error detail /
# timestamps, plus widens the Error Detail column.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/62362015cbd27729e6f576ce078eb82d4e7a1c44/e2e/b3b653f4-b0ac-4d24-bf69-21d55057a0d5.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b84f6477bdc8ae95b5b37d31c8014ec901274b5f/e2e/b3b653f4-b0ac-4d24-bf69-21d55057a0d5.md."
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b84f6477bdc8ae95b5b37d31c8014ec901274b5f/e2e/b3b653f4-b0ac-4d24-bf69-21d55057a0d5.md"
$display = "b3b653f4-b0ac-4d24-bf69-21d55057a0d5.md"

function Update-HandbackRow($ws, $latestHandbackFile, $latestHandbackDateTime) {
    # Widen the "Error Detail" column (P, the 16th column) to width 40.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # I6 = Latest Target File -> becomes a hyperlink to the source .md file
    $ws.Range("I6").Value = $display
    $hl = $ws.Hyperlinks
    $hl.Add($ws.Range("I6"), $targetUrl, "", "", $display) | Out-Null
    # Re-apply the workbook's existing custom HyperLink look (underline +
    # FF6495ED) instead of the engine's default themed hyperlink style.
    $ws.Range("I6").Font.Name = "Calibri"
    $ws.Range("I6").Font.Underline = 2
    $ws.Range("I6").Font.Color = 15570276

    # J6 = Latest Handback File
    $ws.Range("J6").Value = $latestHandbackFile

    # K6 = Latest Handback DateTime
    $ws.Range("K6").Value = $latestHandbackDateTime

    # P6 = Error Detail
    $ws.Range("P6").Value = $errorDetail
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackRow $wsZhCn "b3b653f4-b0ac-4d24-bf69-21d55057a0d5.5394841962e0ab2d249cde9529531a80ad1777b3.zh-cn.xlf" "2016-10-26 07:44:31"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackRow $wsDeDe "b3b653f4-b0ac-4d24-bf69-21d55057a0d5.5394841962e0ab2d249cde9529531a80ad1777b3.de-de.xlf" "2016-10-26 07:44:48"
